# Update dataset reference labels: Breimann23x -> Breimann24x
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Breimann23b" -> "Breimann24b" (cells D2 and D4 share this string)
$ws.Range("D2").Value = "Breimann24b"
$ws.Range("D4").Value = "Breimann24b"

# "Breimann23a" -> "Breimann24a" (cells D5, D6 and D7 share this string)
$ws.Range("D5").Value = "Breimann24a"
$ws.Range("D6").Value = "Breimann24a"
$ws.Range("D7").Value = "Breimann24a"

# Update the sheet view: zoom to 120% and move the selection to D9
$excel.ActiveWindow.Zoom = 120
$ws.Range("D9").Select()
